$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new row at the very top of the data block (new row 2) for
#    "Normalweight concrete, 3000 psi" (single-family home data).
# ---------------------------------------------------------------------------
$ws.Rows("2:2").Insert()
$ws.Range("A2").Value = "Normalweight concrete, 3000 psi"
$ws.Range("B2").Value = 0.05
$ws.Range("C2").Value = 0.05

# Row-insert above inherits formatting from the row above (the header row),
# giving A2 the header style. Restore the "data row" style (same as A3/A4…)
# by copying formats down from a neighboring data cell.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Insert a second new row further down (after "2.5" XPS insulation",
#    before "Glass fiber reinforced concrete (GFRC) Panel") for
#    "Expanded polystyrene (EPS), board".
# ---------------------------------------------------------------------------
$ws.Rows("33:33").Insert()
$ws.Range("A33").Value = "Expanded polystyrene (EPS), board"
$ws.Range("B33").Value = 0.15
$ws.Range("C33").Value = 0.075

# ---------------------------------------------------------------------------
# 3) Re-apply the AutoFilter now so its range covers exactly the (52-row)
#    data block built so far, A1:C53 - done *before* the rows below are
#    appended so the filter doesn't balloon out to the full used range.
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:C53").AutoFilter()

# ---------------------------------------------------------------------------
# 4) Append eight brand-new rows below the existing table (rows 53-60) for
#    additional single-family-home building materials.
# ---------------------------------------------------------------------------
$ws.Range("A53").Value = "Coarse aggregate"
$ws.Range("B53").Value = 0.15
$ws.Range("C53").Value = 0.075

$ws.Range("A54").Value = "Composite wood I-joist"
$ws.Range("B54").Value = 0.1
$ws.Range("C54").Value = 0.05

$ws.Range("A55").Value = "Oriented strandboard (OSB)"
$ws.Range("B55").Value = 0.1
$ws.Range("C55").Value = 0.05

$ws.Range("A56").Value = "Light wood framing"
$ws.Range("B56").Value = 0.1
$ws.Range("C56").Value = 0.05

$ws.Range("A57").Value = "Fiberglass blanket insulation, paper faced"
$ws.Range("B57").Value = 0.1
$ws.Range("C57").Value = 0.05

$ws.Range("A58").Value = "Window frame, vinyl, fixed"
$ws.Range("B58").Value = 0.1
$ws.Range("C58").Value = 0.05

$ws.Range("A59").Value = "SBS modified asphalt shingles"
$ws.Range("B59").Value = 0.1
$ws.Range("C59").Value = 0.05

$ws.Range("A60").Value = "Self-adhering, polymer-modified asphalt sheet underlayment"
$ws.Range("B60").Value = 0.1
$ws.Range("C60").Value = 0.05

# ---------------------------------------------------------------------------
# 5) Update the view: scroll/selection left on B53:C53 (the last appended
#    row), matching where the author finished typing.
# ---------------------------------------------------------------------------
$ws.Range("B53:C53").Select()
